$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.702.64'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '3.272.43'
$ws.Range('E3').Value = '  -1.85%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''574.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.86%  '
$ws.Range('D6').Value = '''172.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.91%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '''0.576'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').Value = '3.264.26'
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('E10').Value = '  -5.48%  '
$ws.Range('E11').Value = '  -2.49%  '
$ws.Range('D12').Value = '''44.78'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.94%  '
$ws.Range('E13').Value = '  -1.54%  '
$ws.Range('D14').Value = '''690.00'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.85%  '
$ws.Range('D15').Value = '3.796.43'
$ws.Range('E15').Value = '  -1.80%  '
$ws.Range('D16').Value = '''8.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.71%  '
$ws.Range('D17').Value = '66.807.58'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').Value = '3.269.05'
$ws.Range('E19').Value = '  -1.95%  '
$ws.Range('D20').Value = '''17.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.26%  '
$ws.Range('D21').Value = '''10.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.16%  '
$ws.Range('D22').Value = '''0.880'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('D23').Value = '''16.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.66%  '
$ws.Range('D24').Value = '''5.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.00%  '
$ws.Range('D25').Value = '''99.63'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.90%  '
$ws.Range('E26').Value = '  -4.11%  '
$ws.Range('E27').Value = '  -5.54%  '
$ws.Range('D28').Value = '''33.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.35%  '
$ws.Range('D29').Value = '''9.02'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.41%  '
$ws.Range('D30').Value = '''8.26'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.91%  '
$ws.Range('D31').Value = '''6.61'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.25%  '
$ws.Range('D32').Value = '''579.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.97%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Value = '''10.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.97%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '3.817.75'
$ws.Range('E34').Value = '  -0.85%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '''0.102'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.43%  '
$ws.Range('D37').Value = '''54.99'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.01%  '
$ws.Range('E38').Value = '  -15.83%  '
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('D41').Value = '''2.55'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.25%  '
$ws.Range('D42').Value = '''31.22'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.42%  '
$ws.Range('D43').Value = '0.0₃0659'
$ws.Range('E43').Value = '  -6.37%  '
$ws.Range('D44').Value = '''0.322'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.40%  '
$ws.Range('E45').Value = '  -8.15%  '
$ws.Range('D46').Value = '''0.0399'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.09%  '
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('E48').Value = '  -1.46%  '
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('E50').Value = '  +3.64%  '
$ws.Range('D51').Value = '''129.24'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.83%  '
